$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$teams = @(
  @{Row=2; Team="Ohio State"; Rating=29.4},
  @{Row=3; Team="Indiana"; Rating=29.2},
  @{Row=4; Team="Oregon"; Rating=25.2},
  @{Row=5; Team="Notre Dame"; Rating=24.5},
  @{Row=6; Team="Alabama"; Rating=23.6},
  @{Row=7; Team="Georgia"; Rating=22.5},
  @{Row=8; Team="Texas Tech"; Rating=22.4},
  @{Row=9; Team="Texas A&M"; Rating=22.1},
  @{Row=10; Team="Utah"; Rating=21.8},
  @{Row=11; Team="USC"; Rating=20.100000000000001},
  @{Row=12; Team="Texas"; Rating=19.8},
  @{Row=13; Team="Miami"; Rating=19.600000000000001},
  @{Row=14; Team="Ole Miss"; Rating=18.5},
  @{Row=15; Team="Oklahoma"; Rating=17.7},
  @{Row=16; Team="Tennessee"; Rating=16.2},
  @{Row=17; Team="Michigan"; Rating=15.9},
  @{Row=18; Team="Vanderbilt"; Rating=15.9},
  @{Row=19; Team="Iowa"; Rating=15.7},
  @{Row=20; Team="BYU"; Rating=15},
  @{Row=21; Team="Penn State"; Rating=14.9},
  @{Row=22; Team="LSU"; Rating=14.4},
  @{Row=23; Team="Missouri"; Rating=13.7},
  @{Row=24; Team="Illinois"; Rating=12.5},
  @{Row=25; Team="Auburn"; Rating=12.3},
  @{Row=26; Team="Florida State"; Rating=12.3},
  @{Row=27; Team="Washington"; Rating=12.1},
  @{Row=28; Team="SMU"; Rating=11.1},
  @{Row=29; Team="Nebraska"; Rating=10.6},
  @{Row=30; Team="Clemson"; Rating=10.3},
  @{Row=31; Team="South Florida"; Rating=10.199999999999999},
  @{Row=32; Team="Florida"; Rating=10.199999999999999},
  @{Row=33; Team="Pittsburgh"; Rating=10},
  @{Row=34; Team="TCU"; Rating=9.9},
  @{Row=35; Team="Louisville"; Rating=9.9},
  @{Row=36; Team="Georgia Tech"; Rating=9.8000000000000007},
  @{Row=37; Team="Cincinnati"; Rating=9.4},
  @{Row=38; Team="Iowa State"; Rating=8.9},
  @{Row=39; Team="Kansas State"; Rating=8.8000000000000007},
  @{Row=40; Team="Arkansas"; Rating=8.6999999999999993},
  @{Row=41; Team="Arizona State"; Rating=8.5},
  @{Row=42; Team="Arizona"; Rating=8.4},
  @{Row=43; Team="South Carolina"; Rating=8.4},
  @{Row=44; Team="Kentucky"; Rating=7.9},
  @{Row=45; Team="Mississippi State"; Rating=7.7},
  @{Row=46; Team="Baylor"; Rating=7.3},
  @{Row=47; Team="Duke"; Rating=6.9},
  @{Row=48; Team="James Madison"; Rating=6.2},
  @{Row=49; Team="Kansas"; Rating=6.2},
  @{Row=50; Team="Virginia"; Rating=6.1},
  @{Row=51; Team="North Texas"; Rating=5.7},
  @{Row=52; Team="Memphis"; Rating=5},
  @{Row=53; Team="Boise State"; Rating=4.5},
  @{Row=54; Team="NC State"; Rating=4.2},
  @{Row=55; Team="Houston"; Rating=4.2},
  @{Row=56; Team="East Carolina"; Rating=3.7},
  @{Row=57; Team="Northwestern"; Rating=3.3},
  @{Row=58; Team="UCF"; Rating=3},
  @{Row=59; Team="Tulane"; Rating=2.4},
  @{Row=60; Team="Rutgers"; Rating=2.1},
  @{Row=61; Team="Minnesota"; Rating=1.9},
  @{Row=62; Team="Wisconsin"; Rating=1.4},
  @{Row=63; Team="Toledo"; Rating=1.3},
  @{Row=64; Team="Old Dominion"; Rating=0.6},
  @{Row=65; Team="Wake Forest"; Rating=0.4},
  @{Row=66; Team="Maryland"; Rating=0.3},
  @{Row=67; Team="Michigan State"; Rating=0},
  @{Row=68; Team="UCLA"; Rating=-0.3},
  @{Row=69; Team="Colorado"; Rating=-1.2},
  @{Row=70; Team="Purdue"; Rating=-1.5},
  @{Row=71; Team="UNLV"; Rating=-1.7},
  @{Row=72; Team="Virginia Tech"; Rating=-1.8},
  @{Row=73; Team="San Diego State"; Rating=-2.1},
  @{Row=74; Team="West Virginia"; Rating=-2.2000000000000002},
  @{Row=75; Team="UConn"; Rating=-2.5},
  @{Row=76; Team="UTSA"; Rating=-3.1},
  @{Row=77; Team="California"; Rating=-3.3},
  @{Row=78; Team="Navy"; Rating=-3.3},
  @{Row=79; Team="Washington State"; Rating=-3.4},
  @{Row=80; Team="Ohio"; Rating=-3.9},
  @{Row=81; Team="New Mexico"; Rating=-4.2},
  @{Row=82; Team="Army"; Rating=-4.3},
  @{Row=83; Team="Syracuse"; Rating=-4.5},
  @{Row=84; Team="Southern Miss"; Rating=-4.9000000000000004},
  @{Row=85; Team="Utah State"; Rating=-5.0999999999999996},
  @{Row=86; Team="Louisiana Tech"; Rating=-5.2},
  @{Row=87; Team="Hawai'i"; Rating=-5.3},
  @{Row=88; Team="North Carolina"; Rating=-5.3},
  @{Row=89; Team="Stanford"; Rating=-6.1},
  @{Row=90; Team="Fresno State"; Rating=-6.7},
  @{Row=91; Team="Texas State"; Rating=-6.7},
  @{Row=92; Team="Miami (OH)"; Rating=-6.8},
  @{Row=93; Team="Kennesaw State"; Rating=-7.3},
  @{Row=94; Team="Troy"; Rating=-7.4},
  @{Row=95; Team="Boston College"; Rating=-7.5},
  @{Row=96; Team="Air Force"; Rating=-7.8},
  @{Row=97; Team="Marshall"; Rating=-8.1999999999999993},
  @{Row=98; Team="Western Michigan"; Rating=-8.3000000000000007},
  @{Row=99; Team="San José State"; Rating=-8.4},
  @{Row=100; Team="Temple"; Rating=-8.4},
  @{Row=101; Team="Wyoming"; Rating=-8.6},
  @{Row=102; Team="Western Kentucky"; Rating=-10.1},
  @{Row=103; Team="Liberty"; Rating=-10.199999999999999},
  @{Row=104; Team="Oregon State"; Rating=-10.5},
  @{Row=105; Team="Jacksonville State"; Rating=-10.6},
  @{Row=106; Team="Georgia Southern"; Rating=-11},
  @{Row=107; Team="Louisiana"; Rating=-11},
  @{Row=108; Team="Arkansas State"; Rating=-11.1},
  @{Row=109; Team="South Alabama"; Rating=-11.1},
  @{Row=110; Team="App State"; Rating=-11.3},
  @{Row=111; Team="Buffalo"; Rating=-11.4},
  @{Row=112; Team="Central Michigan"; Rating=-11.8},
  @{Row=113; Team="Florida Atlantic"; Rating=-12.6},
  @{Row=114; Team="Colorado State"; Rating=-13.1},
  @{Row=115; Team="Missouri State"; Rating=-13.3},
  @{Row=116; Team="Coastal Carolina"; Rating=-13.3},
  @{Row=117; Team="Bowling Green"; Rating=-14.3},
  @{Row=118; Team="Delaware"; Rating=-14.3},
  @{Row=119; Team="Florida International"; Rating=-14.7},
  @{Row=120; Team="UTEP"; Rating=-14.9},
  @{Row=121; Team="Oklahoma State"; Rating=-15.3},
  @{Row=122; Team="UAB"; Rating=-15.8},
  @{Row=123; Team="Rice"; Rating=-16.100000000000001},
  @{Row=124; Team="Eastern Michigan"; Rating=-16.2},
  @{Row=125; Team="Northern Illinois"; Rating=-16.5},
  @{Row=126; Team="Tulsa"; Rating=-17.3},
  @{Row=127; Team="New Mexico State"; Rating=-17.899999999999999},
  @{Row=128; Team="Akron"; Rating=-19.100000000000001},
  @{Row=129; Team="Nevada"; Rating=-19.7},
  @{Row=130; Team="Georgia State"; Rating=-20},
  @{Row=131; Team="UL Monroe"; Rating=-20.100000000000001},
  @{Row=132; Team="Ball State"; Rating=-20.2},
  @{Row=133; Team="Middle Tennessee"; Rating=-21.1},
  @{Row=134; Team="Kent State"; Rating=-22.8},
  @{Row=135; Team="Sam Houston"; Rating=-23},
  @{Row=136; Team="Charlotte"; Rating=-23.9},
  @{Row=137; Team="Massachusetts"; Rating=-30}
)

foreach ($item in $teams) {
  $ws.Cells.Item($item.Row, 1).Value = $item.Team
  $ws.Cells.Item($item.Row, 2).Value = $item.Rating
}

$ws.Range("C136").Select()